$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old header row (row 6), pushing it to row 7
$ws.Rows.Item(6).Insert()

# Write the new deck-name row in column A
$ws.Cells.Item(6, 1).Value2 = "#deck:(KR) eohwi"

# Match the formatting (font/size) used by the other "#..." meta rows (A2:A5)
$ws.Cells.Item(6, 1).Font.Name = $ws.Cells.Item(2, 1).Font.Name
$ws.Cells.Item(6, 1).Font.Size = $ws.Cells.Item(2, 1).Font.Size
$ws.Cells.Item(6, 1).Font.Bold = $ws.Cells.Item(2, 1).Font.Bold

# The meta rows (2-6) only ever use column A - clear any stray formatting that
# spilled into B:D on insert so the row only contains the one cell
$ws.Range("B2:D6").Clear()

# Keep all of the meta rows (2-6) the same height, consistent with the rest
$ws.Range("A5:A6").EntireRow.RowHeight = 19.5
